$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "303.58" or
# "43.110.57" are not auto-converted to numbers by Excel's smart input parsing.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.110.57'
$ws.Range("E2").Value = '  +2.18%  '
$ws.Range("D3").Value = '2.314.07'
$ws.Range("E3").Value = '  +1.97%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '303.58'
$ws.Range("E5").Value = '  +1.85%  '
$ws.Range("D6").Value = '101.40'
$ws.Range("E6").Value = '  +6.67%  '
$ws.Range("E7").Value = '  +2.89%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +3.91%  '
$ws.Range("D10").Value = '34.92'
$ws.Range("E10").Value = '  +4.72%  '
$ws.Range("E11").Value = '  +1.17%  '
$ws.Range("E12").Value = '  +4.17%  '
$ws.Range("D13").Value = '18.14'
$ws.Range("E13").Value = '  +14.89%  '
$ws.Range("E14").Value = '  +3.32%  '
$ws.Range("D15").Value = '2.691.50'
$ws.Range("D16").Value = '2.355.17'
$ws.Range("E16").Value = '  +3.48%  '
$ws.Range("D17").Value = '0.821'
$ws.Range("E17").Value = '  +5.27%  '
$ws.Range("D18").Value = '43.054.29'
$ws.Range("E18").Value = '  +2.15%  '
$ws.Range("D19").Value = '12.54'
$ws.Range("E19").Value = '  +7.88%  '
$ws.Range("E20").Value = '  +3.15%  '
$ws.Range("D21").Value = '0.0₃0906'
$ws.Range("E21").Value = '  +1.87%  '
$ws.Range("D22").Value = '67.91'
$ws.Range("E22").Value = '  +2.01%  '
$ws.Range("D23").Value = '237.43'
$ws.Range("E23").Value = '  +1.78%  '
$ws.Range("E24").Value = '  +12.93%  '
$ws.Range("E25").Value = '  +1.17%  '
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("D27").Value = '24.86'
$ws.Range("E27").Value = '  +4.17%  '
$ws.Range("D28").Value = '168.38'
$ws.Range("E28").Value = '  +0.46%  '
$ws.Range("D29").Value = '2.10'
$ws.Range("E29").Value = '  -3.78%  '
$ws.Range("D30").Value = '34.31'
$ws.Range("E30").Value = '  +0.77%  '
$ws.Range("D31").Value = '9.21'
$ws.Range("E31").Value = '  +1.25%  '
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("E33").Value = '  +2.93%  '
$ws.Range("D34").Value = '4.67'
$ws.Range("E34").Value = '  +3.27%  '
$ws.Range("D35").Value = '17.20'
$ws.Range("E35").Value = '  +4.63%  '
$ws.Range("E36").Value = '  +4.22%  '
$ws.Range("E37").Value = '  +0.93%  '
$ws.Range("E38").Value = '  +4.09%  '
$ws.Range("D39").Value = '1.80'
$ws.Range("E39").Value = '  +4.55%  '
$ws.Range("E40").Value = '  +1.79%  '
$ws.Range("E41").Value = '  +0.90%  '
$ws.Range("D42").Value = '2.31'
$ws.Range("E42").Value = '  -4.49%  '
$ws.Range("D43").Value = '2.005.92'
$ws.Range("E43").Value = '  +2.23%  '
$ws.Range("E44").Value = '  +3.65%  '
$ws.Range("D45").Value = '10.22'
$ws.Range("E45").Value = '  +6.91%  '
$ws.Range("D46").Value = '17.68'
$ws.Range("E46").Value = '  +1.10%  '
$ws.Range("E47").Value = '  +2.90%  '
$ws.Range("D48").Value = '56.13'
$ws.Range("E48").Value = '  +7.94%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '1.55'
$ws.Range("E49").Value = '  +5.72%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.534.87'
$ws.Range("E50").Value = '  +1.56%  '
$ws.Range("D51").Value = '4.59'
$ws.Range("E51").Value = '  +0.71%  '

# Restore original (default) cell formatting on column D now that the text
# values are safely stored, so no stray number-format styling is left behind.
$ws.Range("D2:D51").ClearFormats()
